$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Row -> (old R2 text, new R2 text). Column 5 is the R2 column.
# (Row numbers follow Word's native Cell() indexing, in which the
# section-heading rows, e.g. "Group vs. solo foraging", each occupy one
# row of their own.)
$changes = @(
    @{Row=3;  Old="0.11"; New="0.03"},
    @{Row=4;  Old="0.10"; New="0.02"},
    @{Row=5;  Old="0.11"; New="0.02"},
    @{Row=8;  Old="0.06"; New="0.05"},
    @{Row=9;  Old="0.06"; New="0.02"},
    @{Row=10; Old="0.05"; New="0.02"},
    @{Row=13; Old="0.05"; New="0.00"},
    @{Row=15; Old="0.05"; New="0.01"},
    @{Row=16; Old="0.05"; New="0.00"}
)

foreach ($change in $changes) {
    $cell = $t.Cell($change.Row, 5)
    $range = $cell.Range
    # Trim trailing cell-mark / paragraph-mark characters so only the
    # visible text is compared/replaced.
    $range.End = $range.End - 1
    if ($range.Text -eq $change.Old) {
        $range.Text = $change.New
    }
}
